$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1,2,3,1,3,6,6,5,1,8,3,5,0,0,3,6,4,2,4,2,4,1,4)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
